$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Wetland number changes from 6 to 7
$ws.Range("B1").Value = 7

# Row 5/6 collapse: "Side"/"11:00-11:27" + "Middle"/"11:31-11:51" rows become
# a single row 5 with "middle" label and a time value, row 6 is cleared.
$ws.Range("B5").Value = "middle"
$ws.Range("C5").Value = 0.51458333333333328
$ws.Range("C5").NumberFormat = "h:mm"

$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Row 8 header row stays the same text, just shared-string index shuffled in
# the file - no content change needed (already Methane/Spot/Vial no/time)

# Data rows 9-11 updated, and a new row 12 appended
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 186
$ws.Range("D9").Value = 0.52222222222222225
$ws.Range("D9").NumberFormat = "h:mm"

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 148
$ws.Range("D10").Value = 0.52430555555555558
$ws.Range("D10").NumberFormat = "h:mm"

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 260
$ws.Range("D11").Value = 0.53194444444444444
$ws.Range("D11").NumberFormat = "h:mm"

$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 179
$ws.Range("D12").Value = 0.53263888888888888
$ws.Range("D12").NumberFormat = "h:mm"

# Selection moves to C10
$ws.Range("C10").Select()
